$d = $word.ActiveDocument

# Locate the work-log table row whose "Activity" cell contains the text
# "Continued working on ICT risk sections" (row for 10/4/24).
$table = $d.Tables.Item(1)
$targetRow = $null
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $row = $table.Rows.Item($i)
    if ($row.Cells.Item(2).Range.Text -like "*Continued working on ICT risk sections*") {
        $targetRow = $row
        break
    }
}

$activityCell = $targetRow.Cells.Item(2)
$hoursCell = $targetRow.Cells.Item(3)

# Rewrite the activity text: "Continued working on ICT risk sections" becomes
# "Continued working on ICT risks section. Added more references." split
# across four runs, matching how the edit was actually made (a word tweak
# plus an appended sentence).
$activityRange = $activityCell.Range
$cellStart = $activityRange.Start

$newText = "Continued working on ICT risks section. Added more references."
$activityRange.Text = $newText

# Force run boundaries at the edit points by touching (and reverting) a
# character formatting property on each new segment - this keeps the
# segments as discrete runs instead of being re-merged into one.
$seg2 = $d.Range($cellStart + 29, $cellStart + 30)
$seg2.Font.Bold = $true
$seg2.Font.Bold = $false

$seg3 = $d.Range($cellStart + 30, $cellStart + 38)
$seg3.Font.Bold = $true
$seg3.Font.Bold = $false

$seg4 = $d.Range($cellStart + 38, $cellStart + 62)
$seg4.Font.Bold = $true
$seg4.Font.Bold = $false

# Update the logged hours for that entry from 3 to 6.
$hoursCell.Range.Text = "6"

Write-Output ("Activity now: " + $activityCell.Range.Text)
Write-Output ("Hours now: " + $hoursCell.Range.Text)
